$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "old" method name variants in column C (Method Name) by appending
# "_Old" to reflect that these test cases now refer to a superseded/old
# implementation (the new methods keep their original, un-suffixed names).
$namesToMarkOld = @(
    "verifyPrimeClasses",
    "verifyPrimeSubjects",
    "searchAndViewContentSchool",
    "searchAndViewContentStudent"
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($namesToMarkOld -contains $current) {
        $cell.Value = "$($current)_Old"
    }
}

# Update the sheet view: scroll position and selection to match where the
# author was last working in the sheet.
$ws.Range("C55").Select()
$excel.ActiveWindow.ScrollRow = 29
